# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker-period table (B16:J27) was regenerated: instead of being
# grouped by worker (each worker's "1608" row followed by their "1607"
# row), it is now grouped by period (every worker's "1607" row first,
# followed by every worker's "1608" row). The underlying data (doc type,
# doc number, name, arrears value, base salary) for each worker/period
# pair is unchanged - only the row order / period grouping changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current ("before") table so we can re-derive every
# worker/period combination regardless of how it is currently sorted.
$rows = @(16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27)

$records = @()
foreach ($r in $rows) {
    $record = @{
        TipoDoc = $ws.Cells.Item($r, 2).Value2
        NumDoc  = $ws.Cells.Item($r, 3).Value2
        Nombre  = $ws.Cells.Item($r, 4).Value2
        Periodo = $ws.Cells.Item($r, 5).Value2
        Valor   = $ws.Cells.Item($r, 6).Value2
        Salario = $ws.Cells.Item($r, 7).Value2
    }
    $records += $record
}

# New order: all "1607" rows (in original worker order), then all
# "1608" rows (in original worker order).
$period1607 = $records | Where-Object { $_.Periodo -eq "1607" }
$period1608 = $records | Where-Object { $_.Periodo -eq "1608" }
$ordered = @()
$ordered += $period1607
$ordered += $period1608

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $rec = $ordered[$i]
    $ws.Cells.Item($r, 2).Value = $rec.TipoDoc
    $ws.Cells.Item($r, 3).Value = $rec.NumDoc
    $ws.Cells.Item($r, 4).Value = $rec.Nombre
    $ws.Cells.Item($r, 5).Value = $rec.Periodo
    $ws.Cells.Item($r, 6).Value = $rec.Valor
    $ws.Cells.Item($r, 7).Value = $rec.Salario
}
